$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.219.61"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "2.437.17"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.66"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.77"
$ws.Range("E6").Value = "  -4.18%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").Value = "2.433.75"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("E10").Value = "  -5.05%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("E13").Value = "  -3.99%  "
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("E15").Value = "  -5.77%  "
$ws.Range("D16").Value = "2.868.68"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").Value = "62.193.01"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "2.432.87"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.94"
$ws.Range("E19").Value = "  -5.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.08"
$ws.Range("E20").Value = "  -4.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.75"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("E23").Value = "  -8.47%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.58"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.25"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "628.57"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "2.564.70"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0948"
$ws.Range("E30").Value = "  -10.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -7.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  -4.97%  "
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.94"
$ws.Range("E35").Value = "  -6.46%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  -7.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.374"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "149.63"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.32"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("E42").Value = "  -4.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.74"
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  -11.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.93"
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("E47").Value = "  -4.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  -4.17%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.56"
$ws.Range("E50").Value = "  -9.52%  "
$ws.Range("D51").Value = "0.0₆0234"
$ws.Range("E51").Value = "  +3.34%  "
